$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175745010375977
$ws.Range("B1").Value = 2.955532789230347
$ws.Range("C1").Value = 2.152588605880737
$ws.Range("D1").Value = 1.40204393863678
$ws.Range("E1").Value = 0.9264914989471436
